$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The long "test steps" text currently lives in D3 (row 3). It needs to move
# down to D10 (row 10), taking its row height (84) and wrap-text formatting
# along with it, while D3 itself becomes blank (still keeping its wrap style).
$stepsText = $ws.Range("D3").Value()

# Write the steps text into its new home, D10, and give it the same
# wrap-text formatting that D3 used to have.
$ws.Range("D10").Value() = $stepsText
$ws.Range("D10").WrapText() = $true

# Clear out the old location, leaving the cell's style/formatting intact.
$ws.Range("D3").ClearContents()

# Row 10 now needs the tall row height that row 3 used to need; row 3 goes
# back to the sheet's default (auto) height.
$ws.Rows.Item(10).RowHeight() = 84
$ws.Rows.Item(3).AutoFit()

# Move the active selection from D6 to D7.
$ws.Range("D7").Select() | Out-Null
